# Auto commit at 2025-11-03  9:02:31.72
# 202510-202510.xlsx — fill in the newly received figures for 四方坪站小计 (row 9)
# and update the running total in row 10, then leave the selection where the
# editor last clicked.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9 (四方坪站小计 / Sifangping station subtotal): 卡费 (N9) and 单片机 (O9)
# figures for the period were entered.
$ws.Range("N9").Value = 2853.82
$ws.Range("O9").Value = 970.6

# Row 10 (总  计 / Grand total): 赛菲姆停车系统 (L10) total updated to the
# corrected figure for the period.
$ws.Range("L10").Value = 9115.71

# Leave the cursor where the editor ended up.
$null = $ws.Range("R13").Select()
